# ---------------------------------------------------------------------------
# cmip6_ipsl_citations.xlsx - "Citations" sheet update
# Adds 12 new citation rows (47-58) with text/links/bibtex, and extends the
# pre-formatted (but empty) citation table down to row 106.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Citations" sheet

# 1) Extend the formatted-but-empty table area (previously rows 1:53) down to
#    row 106 by copying the blank row-53 cell formatting (col A bold style,
#    cols B-E wrap style) across the new rows.
$ws.Range("A53:E53").Copy($ws.Range("A54:E106"))

# 2) Fill in the new citation content. Cells are written in the same order the
#    strings first appear in the target workbook's shared-string table so the
#    regenerated table/indices line up with the source edit.
$v = @'
Hourdin2006
'@
$ws.Range("A47").Value2 = $v

$v = @'
10.1007/s00382-006-0158-0
'@
$ws.Range("B47").Value2 = $v

$v = @'
"@article{
    author = {Hourdin, F. and {Musat}, I. and {Bony}, S. and {Braconnot}, P. and {Codron}, F. and {Dufresne}, J.-L. and {Fairhead}, L. and {Filiberti}, M.-A. and {Friedlingstein}, P. and {Grandpeix}, J.-Y. and {Krinner}, G. and {Levan}, P. and {Li}, Z.-X. and {Lott}, F.},
     title = "{The LMDZ4 general circulation model: climate performance and sensitivity to parametrized physics with emphasis on tropical convection}",
   journal = {Climate Dynamics},
      year = 2006,
    volume = 27,
     pages = {787--813},
       doi = {10.1007/s00382-006-0158-0}
}"
'@
$ws.Range("C47").Value2 = $v

$v = @'
The LMDZ4 general circulation model: climate performance and sensitivity to parametrized physics with emphasis on tropical convection
'@
$ws.Range("E47").Value2 = $v

$v = @'
https://link.springer.com/article/10.1007/s00382-006-0158-0
'@
$ws.Range("D47").Value2 = $v

$v = @'
Zhang2017
'@
$ws.Range("A48").Value2 = $v

$v = @'
Krinner2005
'@
$ws.Range("A49").Value2 = $v

$v = @'
"@article{
  title={A dynamic global vegetation model for studies of the coupled atmosphere-biosphere system},
  author={Krinner, Gerhard and Viovy, Nicolas and de Noblet-Ducoudr{\'e}, Nathalie and Og{\'e}e, J{\'e}r{\^o}me and Polcher, Jan and Friedlingstein, Pierre and Ciais, Philippe and Sitch, Stephen and Prentice, I Colin},
  journal={Global Biogeochemical Cycles},
  volume={19},
  number={1},
  year={2005},
  publisher={Wiley Online Library}
}"
'@
$ws.Range("C49").Value2 = $v

$v = @'
A dynamic global vegetation model for studies of the coupled atmosphere-biosphere system
'@
$ws.Range("E49").Value2 = $v

$v = @'
https://agupubs.onlinelibrary.wiley.com/doi/full/10.1029/2003GB002199
'@
$ws.Range("D49").Value2 = $v

$v = @'
10.1029/2003GB002199
'@
$ws.Range("B49").Value2 = $v

$v = @'
Ducoudre1993
'@
$ws.Range("A50").Value2 = $v

$v = @'
https://journals.ametsoc.org/doi/pdf/10.1175/1520-0442(1993)006%3C0248:SANSOP%3E2.0.CO%3B2
'@
$ws.Range("D50").Value2 = $v

$v = @'
deRosnay1998
'@
$ws.Range("A51").Value2 = $v

$v = @'
https://hal.archives-ouvertes.fr/file/index/docid/330830/filename/hess-2-239-1998.pdf
'@
$ws.Range("D51").Value2 = $v

$v = @'
Modelling root water uptake in a complex land surface scheme coupled to a GCM
'@
$ws.Range("E51").Value2 = $v

$v = @'
"@article{
  title={Modelling root water uptake in a complex land surface scheme coupled to a {GCM}},
  author={De Rosnay, Patricia and Polcher, Jan},
  journal={Hydrology and Earth System Sciences Discussions},
  volume={2},
  number={2/3},
  pages={239--255},
  year={1998}
}"
'@
$ws.Range("C51").Value2 = $v

$v = @'
dOrgeval2008
'@
$ws.Range("A52").Value2 = $v

$v = @'
deRosnay2003
'@
$ws.Range("A53").Value2 = $v

$v = @'
Farquhar1980
'@
$ws.Range("A54").Value2 = $v

$v = @'
Collatz1992
'@
$ws.Range("A55").Value2 = $v

$v = @'
Ball1987
'@
$ws.Range("A56").Value2 = $v

$v = @'
"@incollection{
  title={A model predicting stomatal conductance and its contribution to the control of photosynthesis under different environmental conditions},
  author={Ball, J Timothy and Woodrow, Ian E and Berry, Joseph A},
  booktitle={Progress in photosynthesis research},
  pages={221--224},
  year={1987},
  publisher={Springer}
}"
'@
$ws.Range("C56").Value2 = $v

$v = @'
A model predicting stomatal conductance and its contribution to the control of photosynthesis under different environmental conditions
'@
$ws.Range("E56").Value2 = $v

$v = @'
Zhang2016
'@
$ws.Range("A57").Value2 = $v

$v = @'
deRosnay2000
'@
$ws.Range("A58").Value2 = $v

$v = @'
"@article{
  title={Sensitivity of surface fluxes to the number of layers in the soil model used in {GCM}s},
  author={De Rosnay, P and Bruen, M and Polcher, J},
  journal={Geophysical research letters},
  volume={27},
  number={20},
  pages={3329--3332},
  year={2000},
  publisher={Wiley Online Library}
}"
'@
$ws.Range("C58").Value2 = $v

$v = @'
10.1029/2000GL011574
'@
$ws.Range("B58").Value2 = $v

$v = @'
https://agupubs.onlinelibrary.wiley.com/doi/pdf/10.1029/2000GL011574
'@
$ws.Range("D58").Value2 = $v

$v = @'
Sensitivity of surface fluxes to the number of layers in the soil model used in GCMs
'@
$ws.Range("E58").Value2 = $v

$v = @'
"@article{
  title={{SECHIBA}, a new set of parameterizations of the hydrologic exchanges at the land-atmosphere interface within the {LMD} atmospheric general circulation model},
  author={Ducoudr{\'e}, Nathale I and Laval, Katia and Perrier, Alain},
  journal={Journal of Climate},
  volume={6},
  number={2},
  pages={248--273},
  year={1993}
}"
'@
$ws.Range("C50").Value2 = $v

# 3) Row heights: tall rows where real content was added, uniform 20.1pt for
#    still-empty short rows (47-58), and a uniform 130.5pt for the new blank
#    rows 59-106 (matches the target workbook).
$ws.Rows.Item(47).RowHeight = 154.5
$ws.Rows.Item(48).RowHeight = 20.100000000000001
$ws.Rows.Item(49).RowHeight = 153.75
$ws.Rows.Item(50).RowHeight = 131.25
$ws.Rows.Item(51).RowHeight = 118.5
$ws.Rows.Item(52).RowHeight = 20.100000000000001
$ws.Rows.Item(53).RowHeight = 20.100000000000001
$ws.Rows.Item(54).RowHeight = 20.100000000000001
$ws.Rows.Item(55).RowHeight = 20.100000000000001
$ws.Rows.Item(56).RowHeight = 117.75
$ws.Rows.Item(57).RowHeight = 20.100000000000001
$ws.Rows.Item(58).RowHeight = 130.5
$ws.Range("A59:A106").RowHeight = 130.5

# 4) Restore sheet view / selection state to match the saved workbook
#    (scrolled down to the new rows, with 91:106 selected).
$ws.Activate()
$ws.Range("A91:A106").EntireRow.Select()
$excel.ActiveWindow.ScrollRow = 48